$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range("D2") '29.329.33'
Set-TextCell $ws.Range("E2") '  -0.39%  '
Set-TextCell $ws.Range("D3") '1.860.89'
Set-TextCell $ws.Range("E3") '  -1.19%  '
Set-TextCell $ws.Range("E4") '  -0.06%  '
Set-TextCell $ws.Range("D5") '242.44'
Set-TextCell $ws.Range("E5") '  -0.48%  '
Set-TextCell $ws.Range("D6") '0.7001'
Set-TextCell $ws.Range("E6") '  -2.87%  '
Set-TextCell $ws.Range("D7") '1.000'
Set-TextCell $ws.Range("E7") '  -0.12%  '
Set-TextCell $ws.Range("D8") '0.07877'
Set-TextCell $ws.Range("E8") '  -1.14%  '
Set-TextCell $ws.Range("D9") '0.3128'
Set-TextCell $ws.Range("E9") '  -0.77%  '
Set-TextCell $ws.Range("D10") '24.18'
Set-TextCell $ws.Range("E10") '  -3.34%  '
Set-TextCell $ws.Range("D11") '0.07817'
Set-TextCell $ws.Range("E11") '  -4.12%  '
Set-TextCell $ws.Range("D12") '1.865.68'
Set-TextCell $ws.Range("E12") '  -1.30%  '
Set-TextCell $ws.Range("D13") '5.144'
Set-TextCell $ws.Range("E13") '  -2.07%  '
Set-TextCell $ws.Range("D14") '92.35'
Set-TextCell $ws.Range("E14") '  -2.60%  '
Set-TextCell $ws.Range("D15") '0.6952'
Set-TextCell $ws.Range("E15") '  -2.29%  '
Set-TextCell $ws.Range("D16") '6.518'
Set-TextCell $ws.Range("E16") '  +1.70%  '
Set-TextCell $ws.Range("D17") '0.000008505'
Set-TextCell $ws.Range("E17") '  +0.52%  '
Set-TextCell $ws.Range("D18") '29.327.65'
Set-TextCell $ws.Range("E18") '  -0.42%  '
Set-TextCell $ws.Range("D19") '248.56'
Set-TextCell $ws.Range("E19") '  -2.41%  '
Set-TextCell $ws.Range("D20") '2.117.87'
Set-TextCell $ws.Range("E20") '  -1.65%  '
Set-TextCell $ws.Range("D21") '12.97'
Set-TextCell $ws.Range("E21") '  -2.88%  '
Set-TextCell $ws.Range("D22") '0.9998'
Set-TextCell $ws.Range("E22") '  -0.12%  '
Set-TextCell $ws.Range("D23") '7.563'
Set-TextCell $ws.Range("E23") '  -2.92%  '
Set-TextCell $ws.Range("E24") '  -0.08%  '
Set-TextCell $ws.Range("D25") '0.1534'
Set-TextCell $ws.Range("E25") '  -3.50%  '
Set-TextCell $ws.Range("D26") '160.98'
Set-TextCell $ws.Range("E26") '  -1.00%  '
Set-TextCell $ws.Range("D27") '8.937'
Set-TextCell $ws.Range("E27") '  -1.63%  '
Set-TextCell $ws.Range("D28") '18.68'
Set-TextCell $ws.Range("E28") '  -1.64%  '
Set-TextCell $ws.Range("D29") '1.579'
Set-TextCell $ws.Range("E29") '  +4.70%  '
Set-TextCell $ws.Range("D30") '4.284'
Set-TextCell $ws.Range("E30") '  -3.14%  '
Set-TextCell $ws.Range("D31") '4.250'
Set-TextCell $ws.Range("E31") '  -0.96%  '
Set-TextCell $ws.Range("D32") '1.205'
Set-TextCell $ws.Range("E32") '  -1.59%  '
Set-TextCell $ws.Range("E33") '  -1.66%  '
Set-TextCell $ws.Range("D34") '1.884'
Set-TextCell $ws.Range("E34") '  -3.39%  '
Set-TextCell $ws.Range("D35") '0.7512'
Set-TextCell $ws.Range("E35") '  -0.68%  '
Set-TextCell $ws.Range("E36") '  -0.62%  '
Set-TextCell $ws.Range("D37") '2.697'
Set-TextCell $ws.Range("E37") '  -0.16%  '
Set-TextCell $ws.Range("D38") '0.01861'
Set-TextCell $ws.Range("E38") '  -1.51%  '
Set-TextCell $ws.Range("D39") '1.268.93'
Set-TextCell $ws.Range("E39") '  -0.52%  '
Set-TextCell $ws.Range("D40") '2.744'
Set-TextCell $ws.Range("E40") '  -0.81%  '
Set-TextCell $ws.Range("D41") '0.8979'
Set-TextCell $ws.Range("E41") '  -0.77%  '
Set-TextCell $ws.Range("D42") '110.55'
Set-TextCell $ws.Range("E42") '  -2.35%  '
Set-TextCell $ws.Range("D43") '5.946'
Set-TextCell $ws.Range("D44") '69.80'
Set-TextCell $ws.Range("E44") '  -6.47%  '
Set-TextCell $ws.Range("D45") '1.000'
Set-TextCell $ws.Range("E45") '  -0.11%  '
Set-TextCell $ws.Range("D46") '2.019.21'
Set-TextCell $ws.Range("E46") '  -0.96%  '
Set-TextCell $ws.Range("E47") '  -4.91%  '
Set-TextCell $ws.Range("D48") '9.574'
Set-TextCell $ws.Range("E48") '  +0.36%  '
Set-TextCell $ws.Range("D49") '0.5183'
Set-TextCell $ws.Range("E49") '  -0.22%  '
Set-TextCell $ws.Range("D50") '1.782'
Set-TextCell $ws.Range("E50") '  -1.43%  '
Set-TextCell $ws.Range("D51") '0.4273'
Set-TextCell $ws.Range("E51") '  -2.38%  '
